$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - update column F (想去人数) for rows 2,3,4,8,9,10
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 606
$ws1.Range("F3").Value = 133
$ws1.Range("F4").Value = 33
$ws1.Range("F8").Value = 749
$ws1.Range("F9").Value = 3798
$ws1.Range("F10").Value = 74

# Sheet "全部类型" (sheet4.xml) - update column F (想去人数) for rows 2,3,4,8,9,10
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 606
$ws4.Range("F3").Value = 133
$ws4.Range("F4").Value = 33
$ws4.Range("F8").Value = 750
$ws4.Range("F9").Value = 3798
$ws4.Range("F10").Value = 74
